$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed s_val data (regenerated to filter save games).
# Columns: B=TB, C=d2S, D=K, E=IP, G=sum (F=Win unchanged)
$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    3  = @{ B = 0.04763786555579896; C = 0.04240448674262143; D = 0.8054896365839992; E = 0.496779210170732; G = 1.392311199053152 }
    4  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    5  = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 26.21740644021617;  E = 8.660232485948974; G = 38.00504557982321 }
    6  = @{ B = 0.003994804209775715; C = 0.04240448674262143; D = 0.8054896365839992; E = 8.660232485948974; G = 9.512121413485371 }
    7  = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732; G = 7.524616544037286 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 26.21740644021617;  E = 0.496779210170732; G = 31.61296591696135 }
    9  = @{ B = 0.127881588408715;  C = 1.667794583268128;   D = 26.21740644021617;  E = 8.660232485948974; G = 36.67331509784199 }
    10 = @{ B = 0.6753301551942219; C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732; G = 3.645393585217082 }
    11 = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732; G = 7.524616544037286 }
    12 = @{ B = 0.04763786555579896; C = 0.002777888934908601; D = 0.8054896365839992; E = 0.496779210170732; G = 1.352684601245439 }
    13 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    14 = @{ B = 1.459612070389937;  C = 0.3127903958511391;  D = 0.1575252929769615; E = 0.496779210170732; G = 2.42670696938877 }
    15 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    16 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    17 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
